$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-number-looking Price cells to remain as text (matching source formatting)
$textCells = @('D5', 'D6', 'D8', 'D9', 'D12', 'D14', 'D18', 'D20', 'D21', 'D22', 'D23', 'D24', 'D25', 'D26', 'D27', 'D28', 'D29', 'D31', 'D32', 'D33', 'D34', 'D38', 'D40', 'D41', 'D42', 'D45', 'D46', 'D47', 'D48', 'D50')
foreach ($cellref in $textCells) {
    $ws.Range($cellref).NumberFormat = "@"
}

# Apply updated values from the latest cryptos data refresh
$ws.Range('D2').Value = '69.012.73'
$ws.Range('E2').Value = '  +1.70%  '
$ws.Range('D3').Value = '3.764.33'
$ws.Range('E3').Value = '  -0.71%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '623.07'
$ws.Range('E5').Value = '  +3.81%  '
$ws.Range('D6').Value = '165.32'
$ws.Range('E6').Value = '  +1.32%  '
$ws.Range('D7').Value = '3.760.68'
$ws.Range('E7').Value = '  -0.71%  '
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').Value = '0.520'
$ws.Range('E9').Value = '  +1.36%  '
$ws.Range('E10').Value = '  +1.20%  '
$ws.Range('E11').Value = '  +3.10%  '
$ws.Range('D12').Value = '6.72'
$ws.Range('E12').Value = '  +0.92%  '
$ws.Range('E13').Value = '  +0.24%  '
$ws.Range('D14').Value = '35.60'
$ws.Range('E14').Value = '  +1.24%  '
$ws.Range('D15').Value = '4.407.97'
$ws.Range('E15').Value = '  -0.47%  '
$ws.Range('D16').Value = '3.769.08'
$ws.Range('E16').Value = '  -0.40%  '
$ws.Range('D17').Value = '69.033.20'
$ws.Range('E17').Value = '  +1.73%  '
$ws.Range('D18').Value = '17.61'
$ws.Range('E18').Value = '  -2.98%  '
$ws.Range('E19').Value = '  -1.21%  '
$ws.Range('D20').Value = '7.03'
$ws.Range('E20').Value = '  +0.25%  '
$ws.Range('D21').Value = '466.21'
$ws.Range('E21').Value = '  +1.58%  '
$ws.Range('D22').Value = '9.54'
$ws.Range('E22').Value = '  +1.01%  '
$ws.Range('D23').Value = '0.705'
$ws.Range('E23').Value = '  +1.91%  '
$ws.Range('D24').Value = '0.0000146'
$ws.Range('E24').Value = '  +2.73%  '
$ws.Range('D25').Value = '83.14'
$ws.Range('E25').Value = '  +0.19%  '
$ws.Range('D26').Value = '12.02'
$ws.Range('E26').Value = '  +1.25%  '
$ws.Range('D27').Value = '2.15'
$ws.Range('E27').Value = '  +3.50%  '
$ws.Range('D28').Value = '10.02'
$ws.Range('E28').Value = '  +1.32%  '
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  -0.17%  '
$ws.Range('D30').Value = '3.915.81'
$ws.Range('E30').Value = '  -0.58%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '2.66'
$ws.Range('E31').Value = '  +2.65%  '
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').Value = '2.24'
$ws.Range('E32').Value = '  +2.47%  '
$ws.Range('D33').Value = '7.15'
$ws.Range('E33').Value = '  -0.97%  '
$ws.Range('D34').Value = '28.67'
$ws.Range('E34').Value = '  -0.95%  '
$ws.Range('E35').Value = '  +16.79%  '
$ws.Range('E36').Value = '  +0.09%  '
$ws.Range('D37').Value = '3.719.68'
$ws.Range('E37').Value = '  -0.57%  '
$ws.Range('D38').Value = '8.93'
$ws.Range('E38').Value = '  +0.11%  '
$ws.Range('E39').Value = '  +1.96%  '
$ws.Range('D40').Value = '3.35'
$ws.Range('E40').Value = '  +5.29%  '
$ws.Range('D41').Value = '5.81'
$ws.Range('E41').Value = '  +0.17%  '
$ws.Range('D42').Value = '0.964'
$ws.Range('E42').Value = '  -1.44%  '
$ws.Range('E43').Value = '  +0.08%  '
$ws.Range('E44').Value = '  -0.07%  '
$ws.Range('D45').Value = '43.25'
$ws.Range('E45').Value = '  -0.83%  '
$ws.Range('B46').Value = 'Monero'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D46').Value = '152.44'
$ws.Range('E46').Value = '  +0.34%  '
$ws.Range('B47').Value = 'TheGraph'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D47').Value = '0.295'
$ws.Range('E47').Value = '  +0.37%  '
$ws.Range('D48').Value = '46.69'
$ws.Range('E48').Value = '  -0.87%  '
$ws.Range('E49').Value = '  +3.55%  '
$ws.Range('D50').Value = '8.38'
$ws.Range('E50').Value = '  +1.24%  '
$ws.Range('E51').Value = '  +0.10%  '
